# Update automatico via Actualizar 05-17-2020 07-13-29
# Appends the next day's COVID tracking row (16-May-2020) to the
# "Condicion_Pacientes" table on sheet "Hoja1".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$lo = $ws.ListObjects.Item(1)

# Grow the table by one row - this extends the table ref (A1:F64 -> A1:F65),
# the autoFilter range and the sheet dimension automatically.
$newRow = $lo.ListRows.Add()

# Copy the formatting (date number format on column A, centered number
# format on the rest) from the previous last row down onto the new one,
# just like dragging the row's formatting down would in the UI.
$ws.Range("A64:F64").Copy()
$ws.Range("A65:F65").PasteSpecial(-4122)

# Fill in the new day's figures.
$ws.Range("A65").Value = 43967
$ws.Range("B65").Value = 531
$ws.Range("C65").Value = 105
$ws.Range("D65").Value = 256
$ws.Range("E65").Value = 13
$ws.Range("F65").Value = 20

# Match the author's final cursor position/selection.
$ws.Range("E65").Select() | Out-Null
